$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Near the end of the document, drop the duplicated bold heading
#    paragraph ("Play Cash Coaster slot for free!") by scanning paragraphs
#    from the end (robust against the identical, non-bold H1 at the top).
# ---------------------------------------------------------------------------
$dupeText = "Play Cash Coaster slot for free!"
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $txt = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($p.Range.Bold -eq -1 -and $txt -eq $dupeText) {
        $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# 2) Replace the italic summary paragraph's text with the new image prompt,
#    keeping its italic run formatting intact. (Done before step 3 below so
#    the search text is still unique in the document.)
# ---------------------------------------------------------------------------
$newImagePrompt = "Create a feature image for Cash Coaster that captures the fun, upbeat roller coaster theme of the game. The image should be in cartoon style and feature a happy Maya warrior with glasses, who represents the excitement and thrill of the amusement park. This warrior should be shown riding a roller coaster with a big smile on their face, while holding some of the classic amusement park treats like pretzels, cotton candy, and caramel apples. The background could include the roller coaster and the bright neon lights of the Cash Coaster logo, as well as other carnival attractions like a Ferris wheel or a carousel. Overall, the image should convey the playful and exciting vibe of Cash Coaster and entice players to take a ride on this thrilling slot game."

$d.Content.Find.Execute("Our expert review of Cash Coaster slot. Play for free and enjoy one of the top theme park-themed slots.", $true, $false, $false, $false, $false, $true, 1, $false, $newImagePrompt, 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Insert a new "Meta description: ..." paragraph right after the first
#    heading ("Play Cash Coaster slot for free!").
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

$boldText = "Meta description"
$restText = ": Our expert review of Cash Coaster slot. Play for free and enjoy one of the top theme park-themed slots."

$metaStart = $metaPara.Range.Start
$metaRange = $metaPara.Range
$metaRange.Text = $boldText + $restText

$boldRange = $d.Range($metaStart, $metaStart + $boldText.Length)
$boldRange.Font.Bold = 1
